$d = $word.ActiveDocument

# Locate the title paragraph ("Compte-rendu du 04/04/2021") and position
# a collapsed range right after its final character.
$titlePara = $d.Paragraphs(1)
$r = $titlePara.Range
$r.Collapse(0)

# Insert a blank "Sans interligne" paragraph right after the title.
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.MoveStart(4, 1)
$r.Style = "Sansinterligne"

# Insert the "Durée : 34min" paragraph after that blank one.
$r.Collapse(0)
$r.InsertParagraphAfter()
$r.Collapse(0)
$r.MoveStart(4, 1)
$r.Style = "Sansinterligne"

$r.Collapse(0)
$r.InsertAfter("Durée :")
$r.Font.Bold = 1
$r.Collapse(0)
$r.InsertAfter(" 34min")
$r.Font.Bold = 0
